# Apply the "Sig. upregulated and downregulated proteins" update:
# - Add a new classification column (G) using a stricter +/-1.5 logFC cut-off,
#   pushing the old summary note out to column H's style slot and down into a
#   new row 27, and adding the single "down" classification (LPA, row 23)
#   that meets the new cut-off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column G (replaces the old "10 upregulated and 15
#     downregulated" summary note that used to live in G1) ---------------
$ws.Range("G1").Value = "up/down regulated using cut off of +/-1.5"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108   # xlCenter

# --- H1 keeps the old summary cell's centered styling but is now blank --
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter

# --- Only LPA (row 23) exceeds the new +/-1.5 cut-off -------------------
$ws.Range("G23").Value = "down"

# --- Move the original "10 upregulated and 15 downregulated" note down
#     under column F, and add the new column's own summary next to it ----
$ws.Range("F27").Value = "10 upregulated and 15 downregulated "
$ws.Range("F27").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G27").Value = "1 downregulated "

# --- Column widths for the widened/new columns ---------------------------
$ws.Columns("F").ColumnWidth = 32.583333333333336
$ws.Columns("G").ColumnWidth = 34.916666666666664
$ws.Columns("H").ColumnWidth = 32.75

# --- Restore the view state (scroll position + active selection) --------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("G31").Select()
